$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.218.68'
$ws.Range("E2").Value = '  -0.49%  '
$ws.Range("D3").Value = '1.648.35'
$ws.Range("E3").Value = '  -0.85%  '
$ws.Range("E4").Value = '  -0.29%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.65%  '
$ws.Range("E6").Value = '  +1.32%  '
$ws.Range("E7").Value = '  -0.32%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.256'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.64%  '
$ws.Range("E9").Value = '  +0.24%  '
$ws.Range("E10").Value = '  +0.90%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0847'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.30%  '
$ws.Range("D12").Value = '1.878.89'
$ws.Range("E12").Value = '  -0.99%  '
$ws.Range("D13").Value = '1.649.96'
$ws.Range("E13").Value = '  -0.72%  '
$ws.Range("E14").Value = '  -1.63%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.539'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.93%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.59'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.73%  '
$ws.Range("D17").Value = '27.210.13'
$ws.Range("E17").Value = '  -0.51%  '
$ws.Range("D18").Value = '0.0₃0741'
$ws.Range("E18").Value = '  +0.64%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '220.00'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.68%  '
$ws.Range("E20").Value = '  -0.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.86'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.67%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.45'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.48'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.81%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.22'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.10'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.28%  '
$ws.Range("E26").Value = '  -0.25%  '
$ws.Range("E27").Value = '  +0.30%  '
$ws.Range("E28").Value = '  -0.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.82'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.24%  '
$ws.Range("E30").Value = '  -1.01%  '
$ws.Range("E31").Value = '  -0.66%  '
$ws.Range("E32").Value = '  -0.72%  '
$ws.Range("E33").Value = '  +1.29%  '
$ws.Range("E34").Value = '  +1.43%  '
$ws.Range("D35").Value = '1.265.93'
$ws.Range("E35").Value = '  +0.13%  '
$ws.Range("E36").Value = '  +0.08%  '
$ws.Range("E37").Value = '  +0.78%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.544'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.52%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.846'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.90%  '
$ws.Range("E40").Value = '  -0.22%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.810'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.35%  '
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("D44").Value = '1.788.75'
$ws.Range("E44").Value = '  -1.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '62.39'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.90%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '92.18'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.29%  '
$ws.Range("E47").Value = '  -0.35%  '
$ws.Range("E48").Value = '  -0.83%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.73'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.76%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0975'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.68%  '
$ws.Range("E51").Value = '  -0.64%  '
